$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match formatting/style of the row above (s="2" -> vertical center + wrap text)
# but leave column B untouched so no B9 cell gets created.
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C8:E8").Copy()
$ws.Range("C9:E9").PasteSpecial(-4122)  # xlPasteFormats

# Add new row 9 data
$ws.Range("A9").Value = 8
$ws.Range("C9").Value = "Brauchen Sie ein Service?"
$ws.Range("D9").Value = "Dringend:Checkbox"
$ws.Range("E9").Value = "Ja"

# Update the view/selection to match final state
$ws.Range("E9").Select()
